# Auto-generated: apply crypto price/volume updates from the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '56.041.11'
$ws.Cells.Item(2, 5).Value = '  -2.71%  '

$ws.Cells.Item(3, 4).Value = '2.362.92'
$ws.Cells.Item(3, 5).Value = '  -3.80%  '

$ws.Cells.Item(4, 5).Value = '  -0.05%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '500.38'
$ws.Cells.Item(5, 5).Value = '  -2.08%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '128.48'
$ws.Cells.Item(6, 5).Value = '  -3.95%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.00'
$ws.Cells.Item(7, 5).Value = '  +0.13%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.544'
$ws.Cells.Item(8, 5).Value = '  -2.22%  '

$ws.Cells.Item(9, 4).Value = '2.367.13'
$ws.Cells.Item(9, 5).Value = '  -3.62%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.0979'
$ws.Cells.Item(10, 5).Value = '  +0.24%  '

$ws.Cells.Item(11, 5).Value = '  +0.30%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '4.81'
$ws.Cells.Item(12, 5).Value = '  +4.37%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.322'
$ws.Cells.Item(13, 5).Value = '  -0.24%  '

$ws.Cells.Item(14, 4).Value = '2.789.44'
$ws.Cells.Item(14, 5).Value = '  -3.52%  '

$ws.Cells.Item(15, 4).Value = '56.032.18'
$ws.Cells.Item(15, 5).Value = '  -2.75%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '21.38'
$ws.Cells.Item(16, 5).Value = '  -2.45%  '

$ws.Cells.Item(17, 5).Value = '  -1.46%  '

$ws.Cells.Item(18, 4).Value = '2.402.52'
$ws.Cells.Item(18, 5).Value = '  -1.45%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '9.99'
$ws.Cells.Item(19, 5).Value = '  -3.08%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '4.03'
$ws.Cells.Item(20, 5).Value = '  -2.31%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '306.31'
$ws.Cells.Item(21, 5).Value = '  -2.61%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.28'
$ws.Cells.Item(22, 5).Value = '  -1.93%  '

$ws.Cells.Item(23, 5).Value = '  -0.19%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '66.11'
$ws.Cells.Item(24, 5).Value = '  +1.33%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '0.997'
$ws.Cells.Item(25, 5).Value = '  -0.04%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '0.368'
$ws.Cells.Item(26, 5).Value = '  -3.34%  '

$ws.Cells.Item(27, 5).Value = '  -6.11%  '

$ws.Cells.Item(28, 5).Value = '  -4.77%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '171.48'
$ws.Cells.Item(29, 5).Value = '  -0.90%  '

$ws.Cells.Item(30, 4).Value = '0.0₃0710'
$ws.Cells.Item(30, 5).Value = '  -3.31%  '

$ws.Cells.Item(31, 5).Value = '  -3.55%  '

$ws.Cells.Item(33, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.998'
$ws.Cells.Item(33, 5).Value = '  +0.32%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '5.72'
$ws.Cells.Item(34, 5).Value = '  -7.22%  '

$ws.Cells.Item(35, 2).Value = 'Fetch.AI'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.08'
$ws.Cells.Item(35, 5).Value = '  -5.42%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '17.58'
$ws.Cells.Item(36, 5).Value = '  -2.53%  '

$ws.Cells.Item(37, 5).Value = '  -5.87%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '3.73'
$ws.Cells.Item(38, 5).Value = '  -3.45%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '36.05'
$ws.Cells.Item(39, 5).Value = '  -1.78%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.789'
$ws.Cells.Item(40, 5).Value = '  -2.61%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '1.37'
$ws.Cells.Item(41, 5).Value = '  -5.98%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '129.31'
$ws.Cells.Item(42, 5).Value = '  -5.19%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '3.35'
$ws.Cells.Item(43, 5).Value = '  -1.55%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '4.77'
$ws.Cells.Item(44, 5).Value = '  -2.70%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.563'
$ws.Cells.Item(45, 5).Value = '  -2.13%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0901'
$ws.Cells.Item(46, 5).Value = '  -1.78%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '238.83'
$ws.Cells.Item(47, 5).Value = '  -6.92%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.0479'
$ws.Cells.Item(48, 5).Value = '  -2.92%  '

$ws.Cells.Item(49, 5).Value = '  -3.96%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '17.01'
$ws.Cells.Item(50, 5).Value = '  -1.20%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.950'
$ws.Cells.Item(51, 5).Value = '  -0.70%  '
